# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计"),
#    populated with the per-fund holdings table (same shape as the other
#    quarter sheets).
# 2. Prepend a "2022-Q1" row to the "总计" (totals) sheet, pushing the
#    existing rows down by one and renumbering the leading index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell AS TEXT (not auto-coerced to a number)
# by staging it through a scratch cell formatted as Text, then
# Copy/PasteSpecial(values) into the destination. This mirrors how the
# source data stores numeric-looking strings ("9.15", "0.5032", ...) as
# text rather than numbers.
# ---------------------------------------------------------------------
$scratchSheet = $wb.Worksheets.Item(1)
$scratchCell = $scratchSheet.Cells.Item(1048576, 256)
$scratchCell.NumberFormat = "@"

function Set-TextValue($range, [string]$text) {
    $scratchCell.Value = $text
    $scratchCell.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# ---------------------------------------------------------------------
# 1. New "2022-Q1" worksheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1new = $wb.Worksheets.Add($null, $q4)
$q1new.Name = "2022-Q1"

# Header row
$q1new.Cells.Item(1, 2).Value = "基金代码"
$q1new.Cells.Item(1, 3).Value = "基金名称"
$q1new.Cells.Item(1, 4).Value = "基金规模"
$q1new.Cells.Item(1, 5).Value = "股票总仓位"
$q1new.Cells.Item(1, 6).Value = "仓位占比"
$q1new.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1new.Cells.Item(1, 8).Value = "仓位排名"

$fundRows = @(
    @{ A=0; B="257010"; C="国联安小盘精选混合";                         D="9.15";  E="74.36"; F="5.50"; G="0.5032"; H=2 },
    @{ A=1; B="001050"; C="汇添富成长多因子量化策略股票";                 D="11.48"; E="92.68"; F="0.68"; G="0.0781"; H=10 },
    @{ A=2; B="011243"; C="万家惠裕回报6个月持有期混合型证券投资基金A"; D="4.93";  E="23.04"; F="1.33"; G="0.0656"; H=2 },
    @{ A=3; B="002367"; C="国联安安稳灵活配置混合";                      D="2.32";  E="33.99"; F="1.90"; G="0.0441"; H=7 },
    @{ A=4; B="006138"; C="国联安价值优选股票";                          D="0.60";  E="93.30"; F="5.36"; G="0.0322"; H=2 },
    @{ A=5; B="011244"; C="万家惠裕回报6个月持有期混合型证券投资基金C"; D="0.14";  E="23.04"; F="1.33"; G="0.0019"; H=2 }
)

$r = 2
foreach ($row in $fundRows) {
    $q1new.Cells.Item($r, 1).Value = $row.A
    Set-TextValue $q1new.Cells.Item($r, 2) $row.B
    Set-TextValue $q1new.Cells.Item($r, 3) $row.C
    Set-TextValue $q1new.Cells.Item($r, 4) $row.D
    Set-TextValue $q1new.Cells.Item($r, 5) $row.E
    Set-TextValue $q1new.Cells.Item($r, 6) $row.F
    Set-TextValue $q1new.Cells.Item($r, 7) $row.G
    $q1new.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert a "2022-Q1" row at the top of the "总计" totals sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Remember how many data rows existed before the insert (rows 2..lastRow).
$lastRowBefore = $total.Cells.Item($total.Rows.Count, 1).End(-4162).Row   # xlUp

$total.Rows.Item(2).Insert() | Out-Null
$total.Rows.Item(2).ClearFormats() | Out-Null

# Re-apply the same formatting (borders/font/alignment) used by the other
# index-column / data cells to the freshly inserted row.
$total.Range("A3:D3").Copy() | Out-Null
$total.Range("A2:D2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$total.Cells.Item(2, 1).Value = 0
Set-TextValue $total.Cells.Item(2, 2) "2022-Q1"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 0.73

# Renumber the leading index column (0-based) for the rows that shifted
# down, since the source data regenerates this column rather than
# preserving the old numbers.
$newLastRow = $lastRowBefore + 1
for ($row = 3; $row -le $newLastRow; $row++) {
    $total.Cells.Item($row, 1).Value = $row - 2
}

# ---------------------------------------------------------------------
# Cleanup: remove the scratch cell we used for staging text values.
# ---------------------------------------------------------------------
$scratchCell.Clear() | Out-Null
